$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 191 (rows 191-252 shift down to 192-253,
# extending the used range to A1:T253).
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A191").Value = 10
$ws.Range("B191").Value = "Vega Modelo de Temuco"
$ws.Range("C191").Value = "La Araucanía"
$ws.Range("D191").Value = 44524
$ws.Range("E191").Value = 9
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100108
$ws.Range("H191").Value = "Tropicales y subtropicales"
$ws.Range("I191").Value = 100108002
$ws.Range("J191").Value = "Mango"
$ws.Range("K191").Value = "Sin especificar"
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 175
$ws.Range("N191").Value = 7000
$ws.Range("O191").Value = 8000
$ws.Range("P191").Value = 7371
$ws.Range("Q191").Value = "$/bandeja 4 kilos"
$ws.Range("R191").Value = "Perú"
$ws.Range("S191").Value = 1843
$ws.Range("T191").Value = 4
